# Daily attendance processing - 2025-11-18 10:25:18
# Reorders the "Recorded By" (column G) entries so that a leading
# "System" token is moved to the end of the comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -ne $val -and $val.StartsWith("System, ")) {
        $parts = $val.Split(",")
        $trimmedParts = @()
        foreach ($p in $parts) {
            $trimmedParts += $p.Trim()
        }
        # Move the leading "System" token to the end
        $rest = $trimmedParts[1..($trimmedParts.Length - 1)]
        $newParts = $rest + @("System")
        $newVal = [string]::Join(", ", $newParts)
        $cell.Value2 = $newVal
    }
}
